$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.344.38'
$ws.Range("E2").Value = '  -3.46%  '
$ws.Range("D3").Value = '3.150.16'
$ws.Range("E3").Value = '  -2.85%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.38%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.148.11'
$ws.Range("E8").Value = '  -2.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.08%  '
$ws.Range("E10").Value = '  -6.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.52'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.474'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.97%  '
$ws.Range("D15").Value = '3.664.66'
$ws.Range("E15").Value = '  -2.97%  '
$ws.Range("D16").Value = '64.333.24'
$ws.Range("E16").Value = '  -3.53%  '
$ws.Range("D18").Value = '3.143.50'
$ws.Range("E18").Value = '  -3.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.94'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '479.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.73%  '
$ws.Range("E22").Value = '  -5.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.74'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.29%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  -4.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.113'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -34.94%  '
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.74'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.20%  '
$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.40%  '
$ws.Range("E35").Value = '  -4.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.28'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.98'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.11%  '
$ws.Range("D38").Value = '0.0₃0715'
$ws.Range("E38").Value = '  -10.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '444.38'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -12.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0395'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.119'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.65%  '
$ws.Range("D44").Value = '2.828.94'
$ws.Range("E44").Value = '  -3.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.266'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.69%  '
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.25%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.114'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '117.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.00%  '
